$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3217
$ws1.Range("F5").Value = 6882
$ws1.Range("F6").Value = 2071
$ws1.Range("F7").Value = 26
$ws1.Range("F13").Value = 151
$ws1.Range("F15").Value = 36

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3217
$ws4.Range("F6").Value = 6882
$ws4.Range("F7").Value = 2071
$ws4.Range("F8").Value = 26
$ws4.Range("F14").Value = 151
$ws4.Range("F16").Value = 36
